$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for_participant")

# Task 2. Perform online tasks ... -> add "per month"
$ws.Range("B8").Value = "Task 2. Perform online tasks that take just about 15 minutes per month for 3 years"

# Step 3. Register & confirm your consent -> add a <a href=/register> link around "Register"/"Registreer"
$ws.Range("B5").Value = "Step 3. <a href=/register>Register</a> & confirm your consent"
$ws.Range("C5").Value = "Stap 3.  <a href=/register>Registreer</a> & bevestig uw toestemming "

# Update the view / selection to match the saved workbook view state
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C5").Select()
